# ISA.xlsx: "changed ISA to be 8 bit instead of 4 bit"
#
# The instruction-format header cells that used to read "Dest (4)",
# "A (4)" and "D (4)" need to become "Dest (8)", "A (8)" and "D (8)"
# respectively. The other header labels (rA, rB, xxxx) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# irmov row: "A (4)" -> "A (8)"
$ws.Range("G2").Value = "A (8)"

# ld / st rows: "D (4)" -> "D (8)"
$ws.Range("G4").Value = "D (8)"
$ws.Range("G5").Value = "D (8)"

# jmp / jz / call rows: "Dest (4)" -> "Dest (8)"
$ws.Range("E6").Value  = "Dest (8)"
$ws.Range("E7").Value  = "Dest (8)"
$ws.Range("E11").Value = "Dest (8)"

# Update the active cell / selection left by the editor.
$ws.Range("J17").Select()
